$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 90; this pushes the existing rows 90-103 down
# to 91-104 (row 104 ends up with the data that used to live in row 103).
$ws.Rows.Item(90).Insert()

# Populate the newly-inserted row 90. It repeats the same record that was
# already in the (now shifted) row 91, except for a new Fecha (column D).
$ws.Cells.Item(90, 1).Value = 9
$ws.Cells.Item(90, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(90, 3).Value = "Metropolitana"
$ws.Cells.Item(90, 4).Value = 44776
$ws.Cells.Item(90, 5).Value = 13
$ws.Cells.Item(90, 6).Value = 100112005
$ws.Cells.Item(90, 7).Value = "Puerro"
$ws.Cells.Item(90, 8).Value = "Sin especificar"
$ws.Cells.Item(90, 9).Value = "Primera"
$ws.Cells.Item(90, 10).Value = 160
$ws.Cells.Item(90, 11).Value = 7000
$ws.Cells.Item(90, 12).Value = 8000
$ws.Cells.Item(90, 13).Value = 7500
$ws.Cells.Item(90, 14).Value = "$/paquete 20 unidades"
$ws.Cells.Item(90, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(90, 16).Value = 375
$ws.Cells.Item(90, 17).Value = 20
$ws.Cells.Item(90, 18).Value = "Hortaliza"
